$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04344787950197155
$ws.Range("H2").Value = 2.765509757537192
$ws.Range("I2").Value = 12.23157152217572
$ws.Range("G3").Value = 0.06655401847966168
$ws.Range("H3").Value = 36.66024898120357
$ws.Range("G4").Value = 0.002254855861345076
$ws.Range("H4").Value = 18.43380044018446
$ws.Range("G5").Value = 0.01013806761129685
$ws.Range("H5").Value = 167.0476944166587
$ws.Range("G6").Value = 0.03500858487396182
$ws.Range("H6").Value = 0.9778180827108
$ws.Range("G7").Value = 0.0693301195796813
$ws.Range("H7").Value = 30.34685627186169
$ws.Range("G8").Value = -0.01599320471687386
$ws.Range("H8").Value = 15.02366562537225
$ws.Range("G9").Value = 0.009920820453143374
$ws.Range("H9").Value = 146.0835387679397
$ws.Range("G10").Value = -0.02703691356061015
$ws.Range("H10").Value = 62.81099136562563
$ws.Range("G11").Value = -0.06522983704713034
$ws.Range("H11").Value = 29.10651160732808
$ws.Range("G12").Value = -0.2056719762074953
$ws.Range("H12").Value = 15.86497493231495
$ws.Range("G13").Value = -0.2774113310202931
$ws.Range("H13").Value = -0.9454094262223687
$ws.Range("G14").Value = -0.05604359314843363
$ws.Range("H14").Value = -51.07232568658776
$ws.Range("G15").Value = 0.004089380227765119
$ws.Range("H15").Value = 111.7607718231379
$ws.Range("G16").Value = 0.1160907157857339
$ws.Range("H16").Value = -7.355841506739551
$ws.Range("G17").Value = 0.1487563148486218
$ws.Range("H17").Value = 6.064737167727823
$ws.Range("G18").Value = 0.1240214352290053
$ws.Range("H18").Value = -0.570250334957354
$ws.Range("G19").Value = 0.1441954999318305
$ws.Range("H19").Value = 8.243630554091713
$ws.Range("G20").Value = 0.0411701382858842
$ws.Range("H20").Value = 19.90331761728163
$ws.Range("G21").Value = 0.04720029833052352
$ws.Range("H21").Value = -18.67606812816062
$ws.Range("G22").Value = -0.08162493923817364
$ws.Range("H22").Value = -2.225422651433719
$ws.Range("G23").Value = -0.09256758143142083
$ws.Range("H23").Value = -47.99817495179491
$ws.Range("G24").Value = 0.1024948367245592
$ws.Range("H24").Value = -13.22417822169135
$ws.Range("G25").Value = 0.1508310828205763
$ws.Range("H25").Value = 19.54554340998632
$ws.Range("G26").Value = 0.05440967085292964
$ws.Range("H26").Value = 9.46715034186661
$ws.Range("G27").Value = 0.07296129396765851
$ws.Range("H27").Value = -15.82351065339667
$ws.Range("G28").Value = -0.07848184000852908
$ws.Range("H28").Value = -23.41069122937963
$ws.Range("G29").Value = -0.09581820175538164
$ws.Range("H29").Value = -34.6286517985291
$ws.Range("G30").Value = 0.06976294889404065
$ws.Range("H30").Value = 9.504002714600897
$ws.Range("G31").Value = 0.05752843184266525
$ws.Range("H31").Value = -5.037976532788671
$ws.Range("G32").Value = 0.09196394665685635
$ws.Range("H32").Value = -6.411339462978738
$ws.Range("G33").Value = 0.1200931092194388
$ws.Range("H33").Value = 45.94827680835316
$ws.Range("G34").Value = -0.003441297669099602
$ws.Range("H34").Value = -113.2077111687932
$ws.Range("G35").Value = 0.02100440318482725
$ws.Range("H35").Value = 287.370197689012
$ws.Range("G36").Value = -0.005268503932310815
$ws.Range("H36").Value = -1078.899036154643
$ws.Range("G37").Value = 0.01675785528476442
$ws.Range("H37").Value = 233.4833649210064
$ws.Range("G38").Value = 0.1164585412408932
$ws.Range("H38").Value = 8.578694972603808
$ws.Range("G39").Value = 0.09028487993872551
$ws.Range("H39").Value = 5.395689778294456
$ws.Range("G40").Value = 0.006005686031698925
$ws.Range("H40").Value = 102.1941350559972
$ws.Range("G41").Value = 0.0389189482726801
$ws.Range("H41").Value = 159.5159453216366
$ws.Range("G42").Value = 0.09170937333015954
$ws.Range("H42").Value = -9.140279112611211
$ws.Range("G43").Value = 0.1177527299362797
$ws.Range("H43").Value = -1.990747156759112
$ws.Range("G44").Value = 0.0491444321719644
$ws.Range("H44").Value = 37.70735485454964
$ws.Range("G45").Value = 0.03798432438238712
$ws.Range("H45").Value = 132.0337218356641
$ws.Range("G46").Value = 0.0391652284918749
$ws.Range("H46").Value = 8.076238813869098
$ws.Range("G47").Value = 0.04941760701476629
$ws.Range("H47").Value = -2.027558461143399
$ws.Range("G48").Value = 0.03303603845852768
$ws.Range("H48").Value = -22.77925930142872
$ws.Range("G49").Value = 0.07215191025783929
$ws.Range("H49").Value = 3.851410141084248
$ws.Range("G50").Value = 0.0186388383901609
$ws.Range("H50").Value = 7.909194216393119
$ws.Range("G51").Value = 0.0210880142485254
$ws.Range("H51").Value = 8.31190148751654
$ws.Range("G52").Value = -0.1204548549961306
$ws.Range("H52").Value = -16.35850568530781
$ws.Range("G53").Value = -0.08117965241862291
$ws.Range("H53").Value = 12.10043293083337
$ws.Range("G54").Value = 0.112619748182977
$ws.Range("H54").Value = 54.01483833380751
$ws.Range("G55").Value = 0.1083051751929798
$ws.Range("H55").Value = 74.82207878464079
$ws.Range("G56").Value = 0.03022152251758366
$ws.Range("H56").Value = -13.62617717182404
$ws.Range("G57").Value = 0.01892788659133184
$ws.Range("H57").Value = 227.8397588382397
$ws.Range("G58").Value = 0.03071302817225358
$ws.Range("H58").Value = 22.80001134833809
$ws.Range("G59").Value = 0.02856241746078764
$ws.Range("H59").Value = 20.62531721899061
$ws.Range("G60").Value = 0.02877040246322089
$ws.Range("H60").Value = -11.31921908054758
$ws.Range("G61").Value = 0.04603085265435178
$ws.Range("H61").Value = 263.6504778983349
$ws.Range("G62").Value = 0.06024207122676271
$ws.Range("H62").Value = -0.1997042005256153
$ws.Range("G63").Value = 0.04009162409106667
$ws.Range("H63").Value = 23.02009014677122
$ws.Range("G64").Value = 0.0193694719599505
$ws.Range("H64").Value = -52.20504010694721
$ws.Range("G65").Value = 0.07039927039404606
$ws.Range("H65").Value = 25.57413115948465
$ws.Range("G66").Value = 0.09429735033204828
$ws.Range("H66").Value = 0.794455514549083
$ws.Range("G67").Value = 0.109912340553851
$ws.Range("H67").Value = -4.793829501053792
$ws.Range("G68").Value = -0.03134772754013255
$ws.Range("H68").Value = 10.05080859887387
$ws.Range("G69").Value = -0.008265007496269555
$ws.Range("H69").Value = 61.05426067931656
$ws.Range("G70").Value = 0.09381686540574508
$ws.Range("H70").Value = 1.274455029290089
$ws.Range("G71").Value = 0.09944053079157736
$ws.Range("H71").Value = 9.025154182360758
$ws.Range("G72").Value = -0.05984021422066699
$ws.Range("H72").Value = -6.703137118063389
$ws.Range("G73").Value = -0.05465935651818908
$ws.Range("H73").Value = 25.89826819745424
$ws.Range("G74").Value = 0.1000923556165869
$ws.Range("H74").Value = 0.1456262847801648
$ws.Range("G75").Value = 0.1201405669038937
$ws.Range("H75").Value = 23.34251049004335
$ws.Range("G76").Value = 0.01164495108233995
$ws.Range("H76").Value = -54.45960662696459
$ws.Range("G77").Value = 0.005882099509934176
$ws.Range("H77").Value = -58.30971172991612
$ws.Range("G78").Value = 0.05773797745381913
$ws.Range("H78").Value = -10.1731480246396
$ws.Range("G79").Value = 0.07054468906241479
$ws.Range("H79").Value = -8.041837893022901
$ws.Range("G80").Value = -0.1659429824125236
$ws.Range("H80").Value = -0.2020308605129298
$ws.Range("G81").Value = -0.1787762596983573
$ws.Range("H81").Value = 14.90573303652248
$ws.Range("G82").Value = 0.1266096833437853
$ws.Range("H82").Value = 10.38782967895368
$ws.Range("G83").Value = 0.2018448060136552
$ws.Range("H83").Value = 13.40729119546635
$ws.Range("G84").Value = 0.04477825283116746
$ws.Range("H84").Value = 87.8467660472149
$ws.Range("G85").Value = 0.0533870805692313
$ws.Range("H85").Value = -13.29879756555182
